$d = $word.ActiveDocument

# Delete the trailing paragraphs (old paragraphs 18-26) that are removed entirely.
$delStart = $d.Paragraphs.Item(18).Range.Start
$delEnd = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()

$d.Paragraphs.Item(1).Range.Text = 'המאמר היומי של מייק - 03.01.25:' + [char]11 + 'A PERCOLATION MODEL OF EMERGENCE: ANALYZING TRANSFORMERS TRAINED ON A FORMAL LANGUAGE'
$d.Paragraphs.Item(2).Range.Text = 'מבוא:' + [char]11 + 'רשתות נוירונים מודרניות, במיוחד מודלי שפה גדולים , מציגות מגוון רחב של יכולות, המאפשרות להן לשמש כמערכות בסיס למגוון יישומים. מאמר זה מציע הגדרה פנומנולוגית של אמרגנטיות בהקשר של רשתות נוירונים, תוך התמקדות באופן שבו מבנים ותהליכים ספציפיים המונחים בבסיס תהליך יצירת דאטה יכולים להוביל לשיפורים פתאומיים בביצועים במשימות ממוקדות יותר.'
$d.Paragraphs.Item(3).Range.Text = 'מושג חשוב:' + [char]11 + 'הפנומנולוגיה היא גישה פילוסופית המתמקדת בחקר מבני התודעה(consciousness) כפי שהם נחווים מנקודת המבט של האדם. היא שואפת לתאר תופעות או הופעת הדברים כפי שהן נתפסות על ידי בני אדם, ללא הנחות מוקדמות או הטיות תיאורטיות. שיטה זו מדגישה את הבנת החוויות כפי שהן נחיות, במטרה לחשוף את המשמעויות הטבועות בהן'
$d.Paragraphs.Item(4).Range.Text = 'יכולות אמרגנטיות(emergent capabilities) ברשתות נוירונים:' + [char]11 + 'החוקרים מגדירים אמרגנטיות ברשתות נוירונים כרכישת מבנים ספציפיים הגורמים לצמיחה פתאומית בביצועים במשימות ספציפיות. הם חוקרים זאת אמפירית באמצעות מערכת ניסויית המבוססת על שפה פורמלית תלוית-הקשר, ומדגימים שטרנספורמרים שאומנו על מחרוזות משפה זו מציגים יכולות אמרגנטיות. ברגע שהמודל לומד את הדקדוק והמבנים הבסיסיים, הביצועים במשימות קשורות משתפרים משמעותית.'
$d.Paragraphs.Item(5).Range.Text = 'הגדרת השפה הפורמלית:' + [char]11 + 'המערכת הניסויית שהוצעה במאמר משתמשת בדקדוק חופשי-הקשר הסתברותי (PCFG) להגדרת שפה פורמלית תלוית-הקשר. הדקדוק כולל:' + [char]11 + [char]11 + 'סימבולים סופיים(terminal symbols): חלקי דיבור הכוללים נושאים, מושאים, פעלים, תארים, פועלים, מילות חיבור ומילות יחס.' + [char]11 + 'סימבולים לא-סופיים: סמלים המגדירים את מבנה המשפטים.' + [char]11 + 'חוקי יצירת טקסט: חוקים המכתיבים כיצד ניתן לשלב סמלים סופיים ולא-סופיים ליצירת משפטים תקפים.'
$d.Paragraphs.Item(6).Range.Text = 'המודל מאומן על משימות כמו יצירה חופשית, פתרון בלבול וייצור מותנה, כאשר מדדי הביצועים נעקבים לאורך תהליך האימון.'
$d.Paragraphs.Item(7).Range.Text = 'משימות ופרוטוקולי הערכת ביצועי מודלים:' + [char]11 + [char]11 + '1. יצירה חופשית של טקסט: המודל מייצר משפטים העומדים בחוקים הדקדוקיים.' + [char]11 + '2. תיקון טקסט לא תקין: המודל מסדר מחדש מחרוזת מבולבלת של מילים ליצירת משפטים תקפים.' + [char]11 + '3. יצירה מותנית: המודל יוצר משפטים על בסיס ישויות או תכונות נתונות.'
$d.Paragraphs.Item(8).Range.Text = 'ההערכה מתבצעת לפי המדדים כוללים בדיקות דקדוקיות, בדיקות טיפוס, דיוק התאמה מדויקת, דיוק פר-טוקן ועוד, המספקים הערכה מקיפה של יכולות המודל.'
$d.Paragraphs.Item(9).Range.Text = 'תוצאות: דינמיקת הלמידה'
$d.Paragraphs.Item(10).Range.Text = 'התוצאות מגלות 3 שלבים מובחנים בדינמיקת הלמידה של המודל:' + [char]11 + [char]11 + '1. שלב ראשוני: המודל לומד מבנים דקדוקיים בסיסיים עם שיפור מינימלי בביצועים.' + [char]11 + '2. ״שינוי פאזה״: מתרחשת עלייה פתאומית בביצועים ברגע שהמודל מתחיל ״להבין את אילוצי שפה״ פשוטים יחסית.' + [char]11 + '3. שלב ההכללה: המודל מדגים ביצועים משופרים במשימות, המעידים על מעבר משינון להכללה.'
$d.Paragraphs.Item(11).Range.Text = 'יכולות אמרגנטיות של מודלים:'
$d.Paragraphs.Item(12).Range.Text = 'החוקרים מבחינים שככל שמודל השפה לומד את הדקדוק ואילוצי הטיפוס, נצפים שיפורי ביצועים משמעותיים במגוון משימות, במיוחד בפתרון בלבול וייצור מותנה. הנוכחות של מבנים ספציפיים מאפשרת למודל לבנות ״שילובים מורכבים ותקינים״ של ישויות ותכונות, המובילים ליכולות אמרגנטיות בתחום השפה.'
$d.Paragraphs.Item(13).Range.Text = 'נקודת מעבר בלמידה:'
$d.Paragraphs.Item(14).Range.Text = 'המאמר דן באופן שבו הופעת יכולות האמרגנטיות קשורה למספר התכונות התיאוריות שהמודל למד. נקודת המעבר, שבה מתרחשים שיפורי ביצועים משמעותיים, קשורה לסקיילינג של תכונות תיאוריות. קביעה זו מאפשרת לחזות מתי יכולות יופיעו ככל שהמודל ממשיך ללמוד.'
$d.Paragraphs.Item(15).Range.Text = 'מסקנה:'
$d.Paragraphs.Item(16).Range.Text = 'מחקר זה תורם להבנת האמרגנטיות ברשתות נוירונים על ידי יצירת מסגרת המגדירה ומאפיינת תכונות אמרגנטיות על בסיס רכישת מבנים בסיסיים על ידי המודל. הממצאים מצביעים על כך שאילוצים דקדוקיים ואילוצי שפה אחרים משמשים כגורמים חשובים בחיזוי התפתחות יכולות במודלים של שפה.'
$d.Paragraphs.Item(17).Range.Text = 'https://arxiv.org/abs/2408.12578'
